# Update NATMI LR-pair TPM values for Dkk4-Kremen1 sheet.
# Base inputs changed:
#   G2:G4 (Ligand average expression value)      : 0.09226200000000001 -> 0.02506566666666667
#   H2:H4 (Ligand total expression value)        : 0.276786             -> 0.075197
#   M2    (Receptor average expression value)    : 3.238669333333334    -> 3.636278
#   N2    (Receptor total expression value)      : 9.716008             -> 10.908834
# All other touched columns (O,P,Q,R,S,T on rows 2-4) are derived values that
# must be recomputed from the above, following the existing NATMI formulas:
#   Q = G * M                                   (edge average expression weight)
#   R = H * N                                   (edge total expression weight)
#   O = M_row / sum(M_rows for same receptor)   (receptor specificity, avg)
#   P = N_row / sum(N_rows for same receptor)   (receptor specificity, total)
#   S = I * O  (I is always 1 in this sheet, so S == O)
#   T = J * P  (J is always 1 in this sheet, so T == P)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update base (input) values ---
$ws.Range("G2").Value = 0.02506566666666667
$ws.Range("G3").Value = 0.02506566666666667
$ws.Range("G4").Value = 0.02506566666666667

$ws.Range("H2").Value = 0.075197
$ws.Range("H3").Value = 0.075197
$ws.Range("H4").Value = 0.075197

$ws.Range("M2").Value = 3.636278
$ws.Range("N2").Value = 10.908834

# --- Recompute derived values for rows 2-4 ---
$lastRow = 4
$sumM = 0.0
$sumN = 0.0
for ($r = 2; $r -le $lastRow; $r++) {
    $sumM += $ws.Cells.Item($r, 13).Value2  # column M
    $sumN += $ws.Cells.Item($r, 14).Value2  # column N
}

for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2   # G
    $h = $ws.Cells.Item($r, 8).Value2   # H
    $i = $ws.Cells.Item($r, 9).Value2   # I
    $j = $ws.Cells.Item($r, 10).Value2  # J
    $m = $ws.Cells.Item($r, 13).Value2  # M
    $n = $ws.Cells.Item($r, 14).Value2  # N

    $o = $m / $sumM
    $p = $n / $sumN
    $q = $g * $m
    $rr = $h * $n
    $s = $i * $o
    $t = $j * $p

    $ws.Cells.Item($r, 15).Value = $o   # O
    $ws.Cells.Item($r, 16).Value = $p   # P
    $ws.Cells.Item($r, 17).Value = $q   # Q
    $ws.Cells.Item($r, 18).Value = $rr  # R
    $ws.Cells.Item($r, 19).Value = $s   # S
    $ws.Cells.Item($r, 20).Value = $t   # T
}
